# Edit: Anapa, Gelen, GK 2021 (added)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Горячий ключ (Hot Key) rows: fill in previously-missing cells ---

# Row 12 (2020): add U12
$ws.Range("U12").Value = 2003
$ws.Range("U12").HorizontalAlignment = -4108

# Row 13 (2021): add D13, E13, N13, P13, Q13, R13
$ws.Range("D13").Value = 8.2
$ws.Range("E13").Value = 274
$ws.Range("N13").Value = 850
$ws.Range("P13").Value = 151.8
$ws.Range("Q13").Value = 76.258
$ws.Range("R13").Value = 1089
$ws.Range("D13:E13").HorizontalAlignment = -4108
$ws.Range("N13").HorizontalAlignment = -4108
$ws.Range("P13:R13").HorizontalAlignment = -4108

# Row 14 (2022): add C14
$ws.Range("C14").Value = 41
$ws.Range("C14").HorizontalAlignment = -4108

# --- Insert two new rows (2020 data) for Anapa and Gelendzhik before row 16 ---
$ws.Rows("16:17").Insert()

# Row 16: Anapa 2020
$ws.Range("A16").Value = "Анапа"
$ws.Range("B16").Value = 2020
$ws.Range("C16").Value = 88.879
$ws.Range("U16").Value = 7130
$ws.Range("D16:T16").Clear()

# Row 17: Gelendzhik 2020
$ws.Range("A17").Value = "Геленджик"
$ws.Range("B17").Value = 2020
$ws.Range("C17").Value = 76.783
$ws.Range("U17").Value = -1278
$ws.Range("D17:T17").Clear()

# Row 18 (shifted from old row 16, Anapa 2021): already has A18, B18, L18 from the
# row shift; fill in the rest of the columns with the newly-added data
$ws.Range("C18").Value = 81.863
$ws.Range("D18").Value = 28.3
$ws.Range("E18").Value = 1017
$ws.Range("F18").Value = 38808
$ws.Range("N18").Value = 3561
$ws.Range("O18").Value = 2466.7
$ws.Range("P18").Value = 289.2
$ws.Range("Q18").Value = 523.164
$ws.Range("R18").Value = 8415
$ws.Range("S18").Value = 36829.6
$ws.Range("T18").Value = 1014.6
$ws.Range("U18").Value = 1556
$ws.Range("C18:F18").HorizontalAlignment = -4108
$ws.Range("N18:U18").HorizontalAlignment = -4108

# Row 19 (shifted from old row 17, Gelendzhik 2021): already has A19, B19, L19 from
# the row shift; fill in the rest of the columns with the newly-added data
$ws.Range("C19").Value = 80.204
$ws.Range("D19").Value = 18.6
$ws.Range("E19").Value = 262
$ws.Range("F19").Value = 42904
$ws.Range("N19").Value = 2565
$ws.Range("O19").Value = 1318.4
$ws.Range("P19").Value = 4092.5
$ws.Range("Q19").Value = 101.312
$ws.Range("R19").Value = 552
$ws.Range("S19").Value = 24283.1
$ws.Range("T19").Value = 992.7
$ws.Range("U19").Value = -258
$ws.Range("C19:F19").HorizontalAlignment = -4108
$ws.Range("N19:U19").HorizontalAlignment = -4108

# Update selection to match the authored workbook state
$ws.Range("D16").Select()
